$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.303.12'
$ws.Range('E2').Value = '  -0.05%  '

$ws.Range('D3').Value = '2.840.73'
$ws.Range('E3').Value = '  +1.80%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '361.05'
$ws.Range('E5').Value = '  +4.14%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '112.78'
$ws.Range('E6').Value = '  -2.68%  '

$ws.Range('E7').Value = '  +4.14%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  +0.02%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.603'
$ws.Range('E9').Value = '  +1.88%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '41.17'
$ws.Range('E10').Value = '  -3.04%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0874'
$ws.Range('E11').Value = '  +1.67%  '

$ws.Range('E12').Value = '  +1.05%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '20.02'
$ws.Range('E13').Value = '  +0.09%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.82'
$ws.Range('E14').Value = '  -1.27%  '

$ws.Range('D15').Value = '3.284.52'
$ws.Range('E15').Value = '  +1.57%  '

$ws.Range('D16').Value = '2.841.10'
$ws.Range('E16').Value = '  +1.07%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.933'
$ws.Range('E17').Value = '  +4.66%  '

$ws.Range('D18').Value = '52.220.27'
$ws.Range('E18').Value = '  +0.07%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.58'
$ws.Range('E19').Value = '  +3.80%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '3.14'
$ws.Range('E20').Value = '  -0.81%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.51'
$ws.Range('E21').Value = '  +1.02%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.0000100'
$ws.Range('E22').Value = '  +2.14%  '

$ws.Range('E23').Value = '  +1.17%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '70.71'
$ws.Range('E24').Value = '  +0.96%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.84'
$ws.Range('E25').Value = '  +3.00%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '27.04'
$ws.Range('E26').Value = '  +0.85%  '

$ws.Range('E27').Value = '  +0.03%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.38'
$ws.Range('E28').Value = '  +1.49%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.26'
$ws.Range('E29').Value = '  +0.19%  '

$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.144'
$ws.Range('E30').Value = '  +2.19%  '

$ws.Range('B31').Value = 'VeChain'
$ws.Range('C31').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0482'
$ws.Range('E31').Value = '  +5.00%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '35.49'
$ws.Range('E32').Value = '  +3.32%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '52.17'
$ws.Range('E33').Value = '  +4.09%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.94'
$ws.Range('E34').Value = '  +3.11%  '

$ws.Range('E35').Value = '  +14.02%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0857'
$ws.Range('E36').Value = '  +2.64%  '

$ws.Range('E37').Value = '  -0.12%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.30'
$ws.Range('E38').Value = '  +2.52%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.06'
$ws.Range('E39').Value = '  -2.49%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '18.51'

$ws.Range('E41').Value = '  +2.15%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '127.34'
$ws.Range('E42').Value = '  +0.77%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.54'
$ws.Range('E43').Value = '  -2.83%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '23.16'
$ws.Range('E44').Value = '  -2.62%  '

$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.38'
$ws.Range('E46').Value = '  +2.13%  '

$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '2.097.31'
$ws.Range('E47').Value = '  +1.98%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.30'
$ws.Range('E48').Value = '  -1.62%  '

$ws.Range('E49').Value = '  +5.91%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.975'
$ws.Range('E50').Value = '  +1.43%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '9.21'
$ws.Range('E51').Value = '  +2.71%  '
